$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header E1: "Python MSE" -> "Rust MSE" (new shared string)
$ws.Range("E1").Value = "Rust MSE"

# C2: Rust R2_score value, styled with Consolas / green font, vertically centered
$rng = $ws.Range("C2")
$rng.Value = 0.65098400000000001
$f = $rng.Font
$f.Name = "Consolas"
$f.Family = 3
$f.Color = 6858364
$rng.VerticalAlignment = -4108

# Propagate the same style (without re-triggering extra font permutations)
# to the other new "Rust" result cells before filling their values in, so
# the new result strings land right after "Rust MSE" in sharedStrings.
$rng.Copy()
$ws.Range("G2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G2").Value = "24.043 ms"
$ws.Range("I2").Value = "1.565 ms"

$ws.Range("E2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E2").Value = 4606896512.2850199

# Restore the original clicked/selected cell noted in the commit
$ws.Range("I33").Select()
